# Refresh the crypto symbol price/volume snapshot (Price = column D,
# Volume(1h) = column E) on Sheet1, matching the Tue Feb 7 17:50:40 UTC 2023
# GitHub Actions data refresh.
#
# These cells are stored as literal text (e.g. "328.63", "0.05%") rather
# than numbers, so each new value is written with a leading apostrophe.
# This tells Excel to keep the entry as text instead of auto-converting a
# numeric-looking string into a floating point number (which would also
# round values like "0.00000000750" and drop their trailing zeros).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'329.23"
$ws.Range("E2").Value = "'0.27%"
$ws.Range("D3").Value = "'44.32"
$ws.Range("E3").Value = "'0.17%"
$ws.Range("D4").Value = "'5.513"
$ws.Range("E4").Value = "'-0.89%"
$ws.Range("D5").Value = "'0.08108"
$ws.Range("E5").Value = "'0.26%"
$ws.Range("D6").Value = "'2.103"
$ws.Range("E6").Value = "'9.70%"
$ws.Range("D7").Value = "'0.9600"
$ws.Range("E7").Value = "'0.78%"
$ws.Range("D8").Value = "'0.1137"
$ws.Range("E8").Value = "'-3.99%"
$ws.Range("D9").Value = "'0.1878"
$ws.Range("E9").Value = "'1.40%"
$ws.Range("D10").Value = "'10.11"
$ws.Range("E10").Value = "'-0.56%"
$ws.Range("D11").Value = "'0.1006"
$ws.Range("E11").Value = "'3.15%"
$ws.Range("D12").Value = "'0.04717"
$ws.Range("E12").Value = "'5.10%"
$ws.Range("E13").Value = "'-0.98%"
$ws.Range("D14").Value = "'0.001257"
$ws.Range("E14").Value = "'-2.26%"
$ws.Range("D15").Value = "'0.04100"
$ws.Range("E15").Value = "'-2.41%"
$ws.Range("D16").Value = "'0.006095"
$ws.Range("E16").Value = "'4.02%"
$ws.Range("D17").Value = "'3.370"
$ws.Range("E17").Value = "'-0.53%"
$ws.Range("D18").Value = "'4.421"
$ws.Range("E18").Value = "'2.50%"
$ws.Range("D19").Value = "'2.621"
$ws.Range("E19").Value = "'2.56%"
$ws.Range("D20").Value = "'0.3308"
$ws.Range("E20").Value = "'-4.46%"
$ws.Range("D21").Value = "'0.1401"
$ws.Range("E21").Value = "'-1.09%"
$ws.Range("D22").Value = "'0.2490"
$ws.Range("E22").Value = "'-0.71%"
$ws.Range("D23").Value = "'0.001306"
$ws.Range("E23").Value = "'4.83%"
$ws.Range("D24").Value = "'0.004354"
$ws.Range("E24").Value = "'-0.21%"
$ws.Range("D25").Value = "'0.0001250"
$ws.Range("E25").Value = "'4.96%"
$ws.Range("D26").Value = "'0.0003739"
$ws.Range("E26").Value = "'-6.02%"
$ws.Range("D38").Value = "'0.02642"
$ws.Range("E38").Value = "'-1.71%"
$ws.Range("D39").Value = "'0.05645"
$ws.Range("E39").Value = "'1.72%"
$ws.Range("D40").Value = "'0.007618"
$ws.Range("E40").Value = "'0.58%"
$ws.Range("D41").Value = "'0.1405"
$ws.Range("E41").Value = "'-0.47%"
$ws.Range("D42").Value = "'0.007360"
$ws.Range("E42").Value = "'-8.29%"
$ws.Range("D43").Value = "'0.001986"
$ws.Range("E43").Value = "'-1.59%"
$ws.Range("D44").Value = "'0.008781"
$ws.Range("E44").Value = "'4.41%"
$ws.Range("D45").Value = "'0.00007099"
$ws.Range("E45").Value = "'-0.80%"
$ws.Range("D46").Value = "'0.00000000750"
$ws.Range("E46").Value = "'-0.09%"
$ws.Range("D47").Value = "'0.0005807"
$ws.Range("E47").Value = "'-0.07%"
$ws.Range("E48").Value = "'54.03%"
$ws.Range("D49").Value = "'0.003361"
$ws.Range("E49").Value = "'-3.40%"
$ws.Range("D50").Value = "'0.00002100"
$ws.Range("E50").Value = "'-0.09%"
$ws.Range("D51").Value = "'0.0002000"
$ws.Range("E51").Value = "'-0.09%"
